$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B" = 0.9533899455498158; "C" = 0.3259424990483808; "E" = 0.4267040067144734; "F" = 0.4443680307746121; "G" = 0.4109261604107246; "H" = 0.4822788916466294; "I" = 0.31806365448832; "N" = 0.7586558342111829 }
    3 = @{ "B" = 0.8334091467887674; "C" = 0.2849640400723956; "E" = 0.3719729806084615; "F" = 0.3878228170618172; "G" = 0.3895634224847981; "H" = 0.4782466569400725; "I" = 0.3181310253427299; "N" = 0.7690097287870827 }
    4 = @{ "B" = 0.7597119079631511; "C" = 0.2597287555325067; "E" = 0.3384851556841113; "F" = 0.3531389305169483; "G" = 0.3769362205838576; "H" = 0.4762017668437863; "I" = 0.3185309183831535; "N" = 0.7757946210837972 }
    5 = @{ "B" = 0.7296721372314323; "C" = 0.2494259333199977; "E" = 0.3248655487302443; "F" = 0.3390132514313251; "G" = 0.3719118374879713; "H" = 0.4754760975210814; "I" = 0.3187833067333798; "N" = 0.778666944632711 }
    6 = @{ "B" = 0.7246835901765962; "C" = 0.2477139798770622; "E" = 0.3226055778975478; "F" = 0.336668177824194; "G" = 0.3710848162166087; "H" = 0.4753620809878925; "I" = 0.318830595740021; "N" = 0.7791503775315221 }
    7 = @{ "B" = 0.7593068121356055; "C" = 0.2595898866708239; "E" = 0.3383013710676437; "F" = 0.3529483938344953; "G" = 0.3768679712327412; "H" = 0.4761915453390628; "I" = 0.3185339610517346; "N" = 0.775832923406746 }
    8 = @{ "B" = 0.9120262241463024; "C" = 0.3118282149498555; "E" = 0.4078068699088391; "F" = 0.4248636149813478; "G" = 0.403457570627026; "H" = 0.4807987034883041; "I" = 0.3180120325488929; "N" = 0.7621370623436476 }
    9 = @{ "B" = 1.211314811976479; "C" = 0.413707161319735; "E" = 0.5451582370615995; "F" = 0.5661985755041457; "G" = 0.4595694294094699; "H" = 0.4932862672492178; "I" = 0.3198666402156221; "N" = 0.7386761574914971 }
    10 = @{ "B" = 1.431152623251364; "C" = 0.4882644680061503; "E" = 0.6468887751503019; "F" = 0.6702781546542269; "G" = 0.5033371460364435; "H" = 0.5046147284691642; "I" = 0.3230308661722958; "N" = 0.7235160879941915 }
    11 = @{ "B" = 1.53116837481997; "C" = 0.5221298938004111; "E" = 0.6933857715410454; "F" = 0.7176906081379002; "G" = 0.5238271074569241; "H" = 0.5102468048190758; "I" = 0.3248718562756352; "N" = 0.7170717744467368 }
    12 = @{ "B" = 1.569044134532078; "C" = 0.5349472630044829; "E" = 0.711027571260459; "F" = 0.7356546913071611; "G" = 0.5316716139712128; "H" = 0.5124492182878555; "I" = 0.3256275713761099; "N" = 0.7146966385333826 }
    13 = @{ "B" = 1.560886816836955; "C" = 0.5321871027758789; "E" = 0.7072265181229511; "F" = 0.7317853510981394; "G" = 0.5299783283257682; "H" = 0.511971777267604; "I" = 0.3254621963402613; "N" = 0.715205264464899 }
    14 = @{ "B" = 1.534284394200313; "C" = 0.5231845178610683; "E" = 0.6948364670813021; "F" = 0.7191683204515869; "G" = 0.5244707564452824; "H" = 0.5104265976436011; "I" = 0.3249328507738127; "N" = 0.7168750631162837 }
    15 = @{ "B" = 1.51798991249774; "C" = 0.5176693162447918; "E" = 0.6872517683494976; "F" = 0.7114413442032514; "G" = 0.5211083921197428; "H" = 0.5094892282412502; "I" = 0.3246162638296539; "N" = 0.7179063575428373 }
    16 = @{ "B" = 1.424616607600001; "C" = 0.4860503070864866; "E" = 0.6438547432163944; "F" = 0.6671810134426437; "G" = 0.5020099189098914; "H" = 0.5042563626028027; "I" = 0.3229187042369475; "N" = 0.7239463509638142 }
    17 = @{ "B" = 1.367338101735243; "C" = 0.4666404198898135; "E" = 0.6172902577191763; "F" = 0.6400460337215605; "G" = 0.4904435628485544; "H" = 0.5011693646611945; "I" = 0.3219807435990987; "N" = 0.727767602622535 }
    18 = @{ "B" = 1.334393971900738; "C" = 0.4554715221086667; "E" = 0.6020315322712548; "F" = 0.6244449056556647; "G" = 0.4838453525832165; "H" = 0.4994388125163027; "I" = 0.3214790009950192; "N" = 0.7300080260488784 }
    19 = @{ "B" = 1.323239805403432; "C" = 0.4516890748801643; "E" = 0.5968686181586946; "F" = 0.6191636801734006; "G" = 0.4816206026700911; "H" = 0.498860582185074; "I" = 0.321315577668102; "N" = 0.7307738960418462 }
    20 = @{ "B" = 1.373435394670594; "C" = 0.4687071324569274; "E" = 0.620115953929627; "F" = 0.642933953830422; "G" = 0.4916691710266718; "H" = 0.5014933167934146; "I" = 0.3220766784301006; "N" = 0.7273564198430762 }
    21 = @{ "B" = 1.542098113377335; "C" = 0.5258289730597312; "E" = 0.6984747688324688; "F" = 0.7228739723492197; "G" = 0.526086129746858; "H" = 0.5108785568960741; "I" = 0.3250867361363987; "N" = 0.7163828321909165 }
    22 = @{ "B" = 1.65234064377313; "C" = 0.5631224488782891; "E" = 0.7498886359077233; "F" = 0.7751780083420101; "G" = 0.549078350218366; "H" = 0.5174188277969449; "I" = 0.3273957754534251; "N" = 0.7095909998831402 }
    23 = @{ "B" = 1.593500898437185; "C" = 0.54322157737613; "E" = 0.7224286587173765; "F" = 0.7472568307915566; "G" = 0.5367606587933551; "H" = 0.5138906922905733; "I" = 0.3261318482902169; "N" = 0.7131810936820742 }
    24 = @{ "B" = 1.370678852040953; "C" = 0.4677728025829992; "E" = 0.6188384147109218; "F" = 0.6416283278902313; "G" = 0.4911149134891843; "H" = 0.5013467204015001; "I" = 0.3220331895336273; "N" = 0.727542179977732 }
    25 = @{ "B" = 1.130363077931634; "C" = 0.3862012304670088; "E" = 0.5078702979538718; "F" = 0.5279251897347308; "G" = 0.4439521252485861; "H" = 0.4895333749679338; "I" = 0.3190519590955603; "N" = 0.7446587381997105 }
}

foreach ($r in $data.Keys) {
    foreach ($c in $data[$r].Keys) {
        $ws.Range("$c$r").Value = $data[$r][$c]
    }
}
